$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Loc"
$ws.Range("B1").Value = "P_max"

# Replace column A values (filenames) with column C values (locations) for rows 2..60
for ($r = 2; $r -le 60; $r++) {
    $loc = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 1).Value = $loc
}

# Delete column C entirely
$ws.Range("C1:C60").EntireColumn.Delete()
